$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21/22: B21 and A22 both become the literal text "0:06:23.5".
# Using a leading apostrophe makes Excel store this as text with the
# "quote prefix" cell flag (same as a user typing '0:06:23.5 into the cell),
# which is what produces the new quotePrefix style + shared string reuse.
$ws.Range("B21").Formula = "'0:06:23.5"
$ws.Range("A22").Formula = "'0:06:23.5"

# Row 30: B30's end-time value shifts to 6.2962962962962964E-3
$ws.Range("B30").Value = 0.0062962962962962964

# Row 31: A31 and B31 become literal text times as well.
# B31 is entered first so the shared-string table order matches the
# target file (0:09:17.5 = index 111, 0:09:04.5 = index 112).
$ws.Range("B31").Formula = "'0:09:17.5"
$ws.Range("A31").Formula = "'0:09:04.5"

# Update the view: scroll so row 10 is at top and select D19
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$null = $ws.Range("D19").Select()
